$wb = $excel.ActiveWorkbook

# Grab references to the three sheets using their CURRENT (pre-edit) names.
$wsAdmin     = $wb.Worksheets.Item("Administrativos")
$wsDocentes  = $wb.Worksheets.Item("Docentes-Conciliadores")
$wsEstud     = $wb.Worksheets.Item("Estudiantes")

# --- 1) Fill the sheet currently named "Administrativos" with two new people ---
$wsAdmin.Range("A2").Value = "Jairo "
$wsAdmin.Range("B2").Value = "Urrego"
$wsAdmin.Range("C2").Value = 52650
$wsAdmin.Range("D2").Value = 31421697
$wsAdmin.Range("E2").Value = "jairo.ug@ugc.edu.co"

$wsAdmin.Range("A3").Value = "Goku"
$wsAdmin.Range("B3").Value = "Saiyajin"
$wsAdmin.Range("C3").Value = 20252285
$wsAdmin.Range("D3").Value = 789524324
$wsAdmin.Range("E3").Value = "goku@ugc.edu.co"

# --- 2) Fill the "Docentes-Conciliadores" sheet with a new teacher/conciliator ---
$wsDocentes.Range("A2").Value = "Vegeta"
$wsDocentes.Range("B2").Value = "Saiyajin"
$wsDocentes.Range("C2").Value = 3546231
$wsDocentes.Range("D2").Value = 3142169745
$wsDocentes.Range("E2").Value = "vegeta.ug@ugc.edu.co"
$wsDocentes.Hyperlinks.Add($wsDocentes.Range("E2"), "mailto:vegeta.ug@ugc.edu.co")
$wsDocentes.Range("F2").Value = 423423

# --- 3) Fill the sheet currently named "Estudiantes" with a new student ---
$wsEstud.Range("A2").Value = "Bulma"
$wsEstud.Range("B2").Value = "Capsula"
$wsEstud.Range("C2").Value = 52650
$wsEstud.Range("D2").Value = 321321654
$wsEstud.Range("E2").Value = "bulma@ugc.edu.co"
$wsEstud.Hyperlinks.Add($wsEstud.Range("E2"), "mailto:bulma@ugc.edu.co")

# --- 4) Swap the "Estudiantes" and "Administrativos" tab names ---
$wsEstud.Name = "Administrativos_tmp"
$wsAdmin.Name = "Estudiantes"
$wsEstud.Name = "Administrativos"

# --- 5) Restore per-sheet selections (cosmetic, matches saved view state) ---
# $wsEstud now carries the "Administrativos" name; $wsAdmin now carries "Estudiantes".
$wsDocentes.Range("C2").Select()
$wsEstud.Range("C2").Select()

# --- 6) Leave the renamed-to-"Estudiantes" sheet (old Administrativos data) active/selected ---
$wsAdmin.Activate()
$wsAdmin.Range("O27").Select()
